$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 234
$ws.Range("I2").Value = 250
$ws.Range("K2").Value = 250
$ws.Range("M2").Value = -137
# Row 15
$ws.Range("H15").Value = 2086.8386
$ws.Range("I15").Value = 2086.8386
$ws.Range("K15").Value = 6260.5158
$ws.Range("M15").Value = -6091.5158
# Row 33
$ws.Range("H33").Value = 1575.1428
$ws.Range("I33").Value = 1575.1428
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1575.1428
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
# Row 95
$ws.Range("H95").Value = 53942
$ws.Range("J95").Value = 53942
$ws.Range("L95").Value = 53942
$ws.Range("N95").Value = -59434
# Row 118
$ws.Range("H118").Value = 787
$ws.Range("I118").Value = 836.5
$ws.Range("K118").Value = 2509.5
$ws.Range("M118").Value = -852.5
# Row 129
$ws.Range("H129").Value = 1792
$ws.Range("I129").Value = 1792
$ws.Range("K129").Value = 5376
$ws.Range("M129").Value = -376
# Row 132
$ws.Range("H132").Value = 2333.3167
$ws.Range("I132").Value = 2245.6667
$ws.Range("J132").Value = 3998.6667
$ws.Range("K132").Value = 6737.000100000001
$ws.Range("L132").Value = 11996.0001
$ws.Range("M132").Value = -4207.000100000001
$ws.Range("N132").Value = -17056.0001

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 7815.7085
$ws.Range("I61").Value = 1838.8
$ws.Range("K61").Value = 1838.8
$ws.Range("M61").Value = -1626.8
# Row 74
$ws.Range("H74").Value = 64217.27
$ws.Range("I74").Value = 101775.81
$ws.Range("K74").Value = 101775.81
$ws.Range("M74").Value = -100901.81
# Row 77
$ws.Range("H77").Value = 64217.27
$ws.Range("I77").Value = 101775.81
$ws.Range("K77").Value = 508879.05
$ws.Range("M77").Value = -504511.05
# Row 98
$ws.Range("H98").Value = 54796
$ws.Range("J98").Value = 54796
$ws.Range("L98").Value = 54796
$ws.Range("N98").Value = -60786
# Row 122
$ws.Range("H122").Value = 16602.732
$ws.Range("I122").Value = 20458.363
$ws.Range("K122").Value = 61375.08900000001
$ws.Range("M122").Value = -58925.08900000001
# Row 132
$ws.Range("H132").Value = 5618.4873
$ws.Range("I132").Value = 2539.818
$ws.Range("J132").Value = 9602.647000000001
$ws.Range("K132").Value = 7619.454000000001
$ws.Range("L132").Value = 28807.941
$ws.Range("M132").Value = -5089.454000000001
$ws.Range("N132").Value = -33867.94100000001
# Row 136
$ws.Range("H136").Value = 7815.7085
$ws.Range("I136").Value = 1838.8
$ws.Range("K136").Value = 5516.4
$ws.Range("M136").Value = -2966.4

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 120
$ws.Range("H120").Value = 59379
$ws.Range("J120").Value = 59379
$ws.Range("L120").Value = 59379
$ws.Range("N120").Value = -69055

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 51973
$ws.Range("J28").Value = 51973
$ws.Range("L28").Value = 51973
$ws.Range("N28").Value = -52463
# Row 31
$ws.Range("H31").Value = 6996.9287
$ws.Range("I31").Value = 2954.8125
$ws.Range("J31").Value = 12386.417
$ws.Range("K31").Value = 2954.8125
$ws.Range("L31").Value = 12386.417
$ws.Range("M31").Value = -2659.8125
$ws.Range("N31").Value = -12976.417
# Row 34
$ws.Range("H34").Value = 6996.9287
$ws.Range("I34").Value = 2954.8125
$ws.Range("J34").Value = 12386.417
$ws.Range("K34").Value = 2954.8125
$ws.Range("L34").Value = 12386.417
$ws.Range("M34").Value = -2752.8125
$ws.Range("N34").Value = -12790.417
# Row 132
$ws.Range("H132").Value = 5105.222
$ws.Range("I132").Value = 1984.1875
$ws.Range("K132").Value = 5952.5625
$ws.Range("M132").Value = -3422.5625

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 103
$ws.Range("H103").Value = 176.25
$ws.Range("J103").Value = 555
$ws.Range("L103").Value = 1665
$ws.Range("N103").Value = -3423
# Row 113
$ws.Range("H113").Value = 2631.1667
$ws.Range("J113").Value = 3374.2666
$ws.Range("L113").Value = 10122.7998
$ws.Range("N113").Value = -14462.7998
# Row 129
$ws.Range("H129").Value = 100976
$ws.Range("I129").Value = 1280.1428
$ws.Range("J129").Value = 333599.66
$ws.Range("K129").Value = 3840.4284
$ws.Range("L129").Value = 1000798.98
$ws.Range("M129").Value = 1159.5716
$ws.Range("N129").Value = -1010798.98

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 1428688.4
$ws.Range("J2").Value = 2857309.8
$ws.Range("L2").Value = 2857309.8
$ws.Range("N2").Value = -2857535.8
# Row 44
$ws.Range("H44").Value = 12115.5
$ws.Range("I44").Value = 10000
$ws.Range("K44").Value = 10000
$ws.Range("M44").Value = -9404
# Row 48
$ws.Range("H48").Value = 7805.6
$ws.Range("I48").Value = 6463.5
$ws.Range("J48").Value = 8700.333000000001
$ws.Range("K48").Value = 6463.5
$ws.Range("L48").Value = 8700.333000000001
$ws.Range("M48").Value = -5978.5
$ws.Range("N48").Value = -9670.333000000001
# Row 102
$ws.Range("H102").Value = 3589.2974
$ws.Range("I102").Value = 3518.0645
$ws.Range("J102").Value = 3957.3333
$ws.Range("K102").Value = 3518.0645
$ws.Range("L102").Value = 3957.3333
$ws.Range("M102").Value = -1896.0645
$ws.Range("N102").Value = -7201.3333
# Row 107
$ws.Range("H107").Value = 1048.2858
$ws.Range("I107").Value = 1237.25
$ws.Range("K107").Value = 1237.25
$ws.Range("M107").Value = 682.75
# Row 122
$ws.Range("H122").Value = 58881.89
$ws.Range("I122").Value = 85456.586
$ws.Range("K122").Value = 256369.758
$ws.Range("M122").Value = -253919.758
# Row 126
$ws.Range("H126").Value = 2598.0667
$ws.Range("I126").Value = 2608.4443
$ws.Range("J126").Value = 2582.5
$ws.Range("K126").Value = 7825.3329
$ws.Range("L126").Value = 7747.5
$ws.Range("M126").Value = -5355.3329
$ws.Range("N126").Value = -12687.5
# Row 132
$ws.Range("H132").Value = 7699.364
$ws.Range("J132").Value = 20000
$ws.Range("L132").Value = 60000
$ws.Range("N132").Value = -65060

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1087.8966
$ws.Range("I16").Value = 1243.3478
$ws.Range("K16").Value = 1243.3478
$ws.Range("M16").Value = -1073.3478
# Row 22
$ws.Range("H22").Value = 11948.857
$ws.Range("I22").Value = 3910.5
$ws.Range("J22").Value = 22666.666
$ws.Range("K22").Value = 3910.5
$ws.Range("L22").Value = 22666.666
$ws.Range("M22").Value = -3615.5
$ws.Range("N22").Value = -23256.666
# Row 27
$ws.Range("H27").Value = 11948.857
$ws.Range("I27").Value = 3910.5
$ws.Range("J27").Value = 22666.666
$ws.Range("K27").Value = 3910.5
$ws.Range("L27").Value = 22666.666
$ws.Range("M27").Value = -3803.5
$ws.Range("N27").Value = -22880.666
# Row 40
$ws.Range("H40").Value = 4120.1577
$ws.Range("I40").Value = 3406.7693
$ws.Range("J40").Value = 5665.8335
$ws.Range("K40").Value = 3406.7693
$ws.Range("L40").Value = 5665.8335
$ws.Range("M40").Value = -3270.7693
$ws.Range("N40").Value = -5937.8335
# Row 46
$ws.Range("H46").Value = 13236703
$ws.Range("I46").Value = 4926875
$ws.Range("K46").Value = 4926875
$ws.Range("M46").Value = -4926687
# Row 68
$ws.Range("H68").Value = 1666.3334
$ws.Range("I68").Value = 999
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 999
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -250
$ws.Range("N68").Value = -3498
# Row 71
$ws.Range("H71").Value = 1666.3334
$ws.Range("I71").Value = 999
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 4995
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -1251
$ws.Range("N71").Value = -17488
# Row 93
$ws.Range("H93").Value = 3653.842
$ws.Range("I93").Value = 4275.8184
$ws.Range("J93").Value = 2798.625
$ws.Range("K93").Value = 4275.8184
$ws.Range("L93").Value = 2798.625
$ws.Range("M93").Value = -3027.8184
$ws.Range("N93").Value = -5294.625
# Row 122
$ws.Range("H122").Value = 3773.2222
$ws.Range("I122").Value = 2845.1738
$ws.Range("K122").Value = 8535.5214
$ws.Range("M122").Value = -6085.5214
# Row 132
$ws.Range("H132").Value = 15159107
$ws.Range("I132").Value = 27782130
$ws.Range("J132").Value = 11480.934
$ws.Range("K132").Value = 83346390
$ws.Range("L132").Value = 34442.802
$ws.Range("M132").Value = -83343860
$ws.Range("N132").Value = -39502.802

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 20017046
$ws.Range("I81").Value = 1491.2
$ws.Range("K81").Value = 2982.4
$ws.Range("M81").Value = -1921.4
# Row 84
$ws.Range("H84").Value = 20017046
$ws.Range("I84").Value = 1491.2
$ws.Range("K84").Value = 14912
$ws.Range("M84").Value = -9608
# Row 113
$ws.Range("H113").Value = 2090
$ws.Range("I113").Value = 1935.1428
$ws.Range("J113").Value = 2188.5454
$ws.Range("K113").Value = 5805.428400000001
$ws.Range("L113").Value = 6565.6362
$ws.Range("M113").Value = -3635.428400000001
$ws.Range("N113").Value = -10905.6362
